# IUserRepository und UserRepositoryDB erweitern (ChangePassword)
# Mark the "ChangePassword" repository task (row 30) as done with a
# completion date, and move the "in Bearbeitung" marker down to the next
# task (row 31, "ChangePassword Methode im UserController erstellen").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Task row 30 is now finished: Zustand -> "done", abgeschlossen am -> date
$ws.Range("C30").Value = "done"
$ws.Range("D30").Value = Get-Date -Year 2019 -Month 3 -Day 7
$ws.Range("D30").NumberFormat = "m/d/yyyy"

# The next task (row 31) is now the one "in Bearbeitung"
$ws.Range("C31").Value = "b"

# Update the active selection to D30, matching where work just happened
$ws.Activate()
$ws.Range("D30").Select()
